$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 47, pushing the existing rows 47-89 down to 48-90.
$ws.Rows(47).Insert()

# Populate the newly inserted row 47 with the new weekly price record.
$ws.Range("A47").Value = 9
$ws.Range("B47").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C47").Value = "Metropolitana"
$ws.Range("D47").Value = 44944
$ws.Range("E47").Value = 13
$ws.Range("F47").Value = 100112029
$ws.Range("G47").Value = "Orégano"
$ws.Range("H47").Value = "Sin especificar"
$ws.Range("I47").Value = "Primera"
$ws.Range("J47").Value = 160
$ws.Range("K47").Value = 20000
$ws.Range("L47").Value = 20000
$ws.Range("M47").Value = 20000
$ws.Range("N47").Value = "`$/docena de atados"
$ws.Range("O47").Value = "Región Metropolitana"
$ws.Range("P47").Value = 6667
$ws.Range("Q47").Value = 3
$ws.Range("R47").Value = "Hortaliza"
